$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'26.756.03"
$ws.Range("E2").Value = "  -2.23%  "

$ws.Range("D3").Value = "'1.798.31"
$ws.Range("E3").Value = "  -1.64%  "

$ws.Range("D4").Value = "'1.000"
$ws.Range("E4").Value = "  -0.03%  "

$ws.Range("D5").Value = "'308.38"
$ws.Range("E5").Value = "  -1.80%  "

$ws.Range("D6").Value = "'1.000"
$ws.Range("E6").Value = "  -0.09%  "

$ws.Range("D7").Value = "'0.4584"
$ws.Range("E7").Value = "  +2.51%  "

$ws.Range("D8").Value = "'0.3716"
$ws.Range("E8").Value = "  -1.37%  "

$ws.Range("D9").Value = "'0.07252"
$ws.Range("E9").Value = "  -3.78%  "

$ws.Range("D10").Value = "'0.8547"
$ws.Range("E10").Value = "  -4.53%  "

$ws.Range("D11").Value = "'20.36"
$ws.Range("E11").Value = "  -3.21%  "

$ws.Range("D12").Value = "'1.808.48"
$ws.Range("E12").Value = "  -0.97%  "

$ws.Range("D13").Value = "'5.309"
$ws.Range("E13").Value = "  -1.69%  "

$ws.Range("D14").Value = "'0.07038"
$ws.Range("E14").Value = "  -1.09%  "

$ws.Range("D15").Value = "'6.483"
$ws.Range("E15").Value = "  -3.93%  "

$ws.Range("D16").Value = "'90.45"
$ws.Range("E16").Value = "  -4.29%  "

$ws.Range("D17").Value = "'1.001"
$ws.Range("E17").Value = "  -0.06%  "

$ws.Range("D18").Value = "'0.000008622"
$ws.Range("E18").Value = "  -2.34%  "

$ws.Range("E19").Value = "  -0.04%  "

$ws.Range("D20").Value = "'14.63"
$ws.Range("E20").Value = "  -3.87%  "

$ws.Range("D21").Value = "'26.763.71"
$ws.Range("E21").Value = "  -2.23%  "

$ws.Range("D22").Value = "'5.288"
$ws.Range("E22").Value = "  +0.07%  "

$ws.Range("D23").Value = "'10.62"
$ws.Range("E23").Value = "  -3.03%  "

$ws.Range("D24").Value = "'2.032.34"
$ws.Range("E24").Value = "  -1.12%  "

$ws.Range("D25").Value = "'1.907"
$ws.Range("E25").Value = "  -4.71%  "

$ws.Range("D26").Value = "'149.50"
$ws.Range("E26").Value = "  -1.36%  "

$ws.Range("B27").Value = "EthereumClassic"
$ws.Range("C27").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D27").Value = "'18.22"
$ws.Range("E27").Value = "  -2.26%  "

$ws.Range("B28").Value = "LidoDAOToken"
$ws.Range("C28").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D28").Value = "'2.145"
$ws.Range("E28").Value = "  -14.20%  "

$ws.Range("D29").Value = "'5.210"
$ws.Range("E29").Value = "  -2.89%  "

$ws.Range("D30").Value = "'114.13"
$ws.Range("E30").Value = "  -3.32%  "

$ws.Range("D31").Value = "'0.08841"
$ws.Range("E31").Value = "  -0.09%  "

$ws.Range("D32").Value = "'0.7530"
$ws.Range("E32").Value = "  -3.62%  "

$ws.Range("D33").Value = "'1.157"
$ws.Range("E33").Value = "  -3.71%  "

$ws.Range("D34").Value = "'4.426"
$ws.Range("E34").Value = "  -1.63%  "

$ws.Range("D35").Value = "'2.884"
$ws.Range("E35").Value = "  -0.28%  "

$ws.Range("E36").Value = "  -0.05%  "

$ws.Range("D37").Value = "'1.114"
$ws.Range("E37").Value = "  +0.46%  "

$ws.Range("D38").Value = "'0.01939"
$ws.Range("E38").Value = "  -2.53%  "

$ws.Range("E39").Value = "  -2.35%  "

$ws.Range("D40").Value = "'2.896"
$ws.Range("E40").Value = "  +0.80%  "

$ws.Range("D41").Value = "'7.153"
$ws.Range("E41").Value = "  -3.19%  "

$ws.Range("D42").Value = "'2.352"
$ws.Range("E42").Value = "  +3.43%  "

$ws.Range("D43").Value = "'0.5216"
$ws.Range("E43").Value = "  -2.15%  "

$ws.Range("D44").Value = "'0.1643"
$ws.Range("E44").Value = "  -4.98%  "

$ws.Range("D45").Value = "'8.477"
$ws.Range("E45").Value = "  -3.92%  "

$ws.Range("D46").Value = "'0.4994"
$ws.Range("E46").Value = "  -3.27%  "

$ws.Range("E47").Value = "  -4.96%  "

$ws.Range("E48").Value = "  -2.10%  "

$ws.Range("E49").Value = "  -0.03%  "

$ws.Range("D50").Value = "'1.642"
$ws.Range("E50").Value = "  -3.84%  "

$ws.Range("D51").Value = "'0.06298"
$ws.Range("E51").Value = "  -1.21%  "
